# Updated cryptos list on Tue Jul 16 18:40:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.692.88'
$ws.Range('E2').Value = '  +1.90%  '

$ws.Range('D3').Value = '3.462.63'
$ws.Range('E3').Value = '  +1.97%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').Value = '''576.54'
$ws.Range('E5').Value = '  -0.11%  '

$ws.Range('D6').Value = '''161.99'
$ws.Range('E6').Value = '  +3.88%  '

$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').Value = '3.466.43'

$ws.Range('E9').Value = '  +8.98%  '

$ws.Range('D10').Value = '''7.35'
$ws.Range('E10').Value = '  -1.77%  '

$ws.Range('E11').Value = '  +4.61%  '

$ws.Range('D12').Value = '''0.443'
$ws.Range('E12').Value = '  +1.49%  '

$ws.Range('D13').Value = '4.062.06'
$ws.Range('E13').Value = '  +2.02%  '

$ws.Range('E14').Value = '  -2.85%  '

$ws.Range('E15').Value = '  +5.84%  '

$ws.Range('D16').Value = '''28.92'
$ws.Range('E16').Value = '  +6.73%  '

$ws.Range('D17').Value = '64.777.74'
$ws.Range('E17').Value = '  +1.86%  '

$ws.Range('D18').Value = '3.445.81'
$ws.Range('E18').Value = '  +1.29%  '

$ws.Range('D19').Value = '''6.39'
$ws.Range('E19').Value = '  +0.21%  '

$ws.Range('D20').Value = '''14.40'
$ws.Range('E20').Value = '  +2.63%  '

$ws.Range('D21').Value = '''391.43'
$ws.Range('E21').Value = '  +0.97%  '

$ws.Range('D22').Value = '''8.20'
$ws.Range('E22').Value = '  -2.77%  '

$ws.Range('D23').Value = '''0.549'
$ws.Range('E23').Value = '  +2.55%  '

$ws.Range('E24').Value = '  +3.22%  '

$ws.Range('E25').Value = '  +0.33%  '

$ws.Range('D26').Value = '''0.0000124'
$ws.Range('E26').Value = '  +20.19%  '

$ws.Range('D27').Value = '''9.47'
$ws.Range('E27').Value = '  -0.48%  '

$ws.Range('E28').Value = '  +0.38%  '

$ws.Range('E29').Value = '  -0.19%  '

$ws.Range('E30').Value = '  +9.77%  '

$ws.Range('E31').Value = '  +7.72%  '

$ws.Range('E32').Value = '  +0.08%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''23.68'
$ws.Range('E33').Value = '  +2.38%  '

$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').Value = '''6.55'
$ws.Range('E34').Value = '  +0.95%  '

$ws.Range('E35').Value = '  +0.12%  '

$ws.Range('D36').Value = '''7.08'
$ws.Range('E36').Value = '  +5.48%  '

$ws.Range('E37').Value = '  +1.66%  '

$ws.Range('D38').Value = '''161.76'
$ws.Range('E38').Value = '  +2.31%  '

$ws.Range('D39').Value = '''1.91'
$ws.Range('E39').Value = '  +1.87%  '

$ws.Range('E40').Value = '  +1.72%  '

$ws.Range('D41').Value = '''27.60'
$ws.Range('E41').Value = '  +0.21%  '

$ws.Range('D42').Value = '2.935.36'
$ws.Range('E42').Value = '  +1.12%  '

$ws.Range('D43').Value = '''4.58'
$ws.Range('E43').Value = '  +6.09%  '

$ws.Range('D44').Value = '''0.0318'
$ws.Range('E44').Value = '  -1.09%  '

$ws.Range('E45').Value = '  +3.46%  '

$ws.Range('D47').Value = '''24.26'
$ws.Range('E47').Value = '  +8.61%  '

$ws.Range('E48').Value = '  +2.91%  '

$ws.Range('D49').Value = '''2.20'
$ws.Range('E49').Value = '  +14.32%  '

$ws.Range('D50').Value = '''0.872'
$ws.Range('E50').Value = '  +6.85%  '

$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '''0.107'
$ws.Range('E51').Value = '  +3.70%  '
